$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.150.75"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "1.836.70"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9989"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.40"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6851"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9994"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3010"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07457"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.11"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.97%  "
$ws.Range("E11").Value = "  -2.03%  "
$ws.Range("D12").Value = "1.837.15"
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.054"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6817"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "87.50"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -6.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.159"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -7.37%  "
$ws.Range("D17").Value = "29.125.19"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008180"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.56%  "
$ws.Range("D19").Value = "2.080.92"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "227.74"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -5.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.52"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.97%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.397"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9994"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1457"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.95"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.753"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.09"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.515"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.270"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.152"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.199"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05171"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7658"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.843"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.134"
$ws.Range("D36").ClearFormats()
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.673"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("D38").Value = "1.306.20"
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01831"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9342"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.798"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "104.62"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.85%  "
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.26"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.48%  "
$ws.Range("D47").Value = "1.981.98"
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5198"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.539"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.771"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05937"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.93%  "
